# 8.10.1 indicator sheet: add a new "2022" data column (S) sourced from the
# existing "2021" column's formatting, revise the 2021 figures that were
# recomputed (R4/R5/R8 - formulas become static reported values), and move
# the active selection to R13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Clone column R's formatting into the new column S (rows 3-8) ---
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial($xlPasteFormats)

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial($xlPasteFormats)

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial($xlPasteFormats)

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial($xlPasteFormats)

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial($xlPasteFormats)

$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Column R (2021) revisions: the a)/b) ratios and adult population were
#     recalculated, so the live formulas are replaced with the reported
#     static values, and the adult-population total is corrected. ---
$ws.Range("R4").Value = 6.9132648934880807
$ws.Range("R5").Value = 42.321589572314856
$ws.Range("R8").Value = 4513063

# --- New column S (2022) values ---
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 6.9031689452913012
$ws.Range("S5").Value = 44.306188104841333
$ws.Range("S6").Value = 318
$ws.Range("S7").Value = 2041
$ws.Range("S8").Value = 4606580

# --- Active cell / selection moves to R13 ---
$ws.Range("R13").Select()
